$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set rows 5 through 24 to a fixed (custom) row height of 14.25 points,
# matching the author's re-upload of the sheet.
$ws.Range("A5:A24").EntireRow.RowHeight = 14.25
